$d = $word.ActiveDocument

$oldText = "Τα διαγράμματα αυτού του αρχείου επιμελήθηκε ο Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "Τα διαγράμματα αυτού του αρχείου επιμελήθηκε ο Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

# Locate the whole credit sentence (it is currently split across several
# differently-formatted runs) and replace it, in one shot, with a single
# plain run containing the updated (2022) link.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false)

if ($found) {
    $rng.Delete()
    $rng.InsertAfter($newText)
} else {
    # Fallback: in case only the year needs swapping (e.g. formatting was
    # already normalised), just update the year inside the URL.
    $d.Content.Find.Execute("GaNight/2018/", $true, $false, $false, $false, $false, $true, 1, $false, "GaNight/2022/", 2) | Out-Null
}

Write-Output "found=$found"
